$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking text (e.g. "27.199.39", "6.90") that must
# remain stored as text, matching the original inlineStr cell type. Temporarily
# force a Text number format while assigning the value, then restore the cells
# original style so no visible formatting change is introduced.
function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '27.199.39'
$ws.Range('E2').Value = '  -0.14%  '
Set-TextValue $ws.Range('D3') '1.630.99'
$ws.Range('E3').Value = '  -0.97%  '
$ws.Range('E4').Value = '  +0.06%  '
Set-TextValue $ws.Range('D5') '216.09'
$ws.Range('E5').Value = '  -0.50%  '
$ws.Range('E6').Value = '  +1.38%  '
$ws.Range('E7').Value = '  +0.12%  '
Set-TextValue $ws.Range('D8') '0.256'
$ws.Range('E8').Value = '  -0.28%  '
$ws.Range('E9').Value = '  -0.97%  '
Set-TextValue $ws.Range('D10') '20.33'
$ws.Range('E10').Value = '  +1.63%  '
Set-TextValue $ws.Range('D11') '0.0849'
$ws.Range('E11').Value = '  +0.69%  '
Set-TextValue $ws.Range('D12') '1.627.10'
$ws.Range('E12').Value = '  -2.36%  '
$ws.Range('E13').Value = '  -0.02%  '
Set-TextValue $ws.Range('D14') '0.543'
$ws.Range('E14').Value = '  +0.42%  '
Set-TextValue $ws.Range('D15') '27.166.52'
$ws.Range('E15').Value = '  -0.20%  '
Set-TextValue $ws.Range('D16') '64.57'
$ws.Range('E16').Value = '  -4.45%  '
$ws.Range('E17').Value = '  -1.00%  '
Set-TextValue $ws.Range('D18') '216.01'
$ws.Range('E18').Value = '  -1.24%  '
Set-TextValue $ws.Range('D20') '6.90'
$ws.Range('E20').Value = '  +0.81%  '
$ws.Range('E21').Value = '  -1.33%  '
Set-TextValue $ws.Range('D22') '2.49'
$ws.Range('E22').Value = '  -0.48%  '
Set-TextValue $ws.Range('D23') '9.09'
$ws.Range('E23').Value = '  -1.16%  '
Set-TextValue $ws.Range('D24') '147.67'
$ws.Range('E24').Value = '  +0.08%  '
$ws.Range('E25').Value = '  +0.13%  '
$ws.Range('E26').Value = '  -3.72%  '
Set-TextValue $ws.Range('D28') '15.59'
$ws.Range('E28').Value = '  -1.45%  '
$ws.Range('E29').Value = '  -0.38%  '
$ws.Range('E30').Value = '  -0.46%  '
Set-TextValue $ws.Range('D31') '3.40'
$ws.Range('E31').Value = '  +0.41%  '
$ws.Range('E32').Value = '  -1.06%  '
Set-TextValue $ws.Range('D33') '1.315.44'
$ws.Range('E33').Value = '  +3.99%  '
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range('D34') '1.56'
$ws.Range('E34').Value = '  -1.66%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range('D35') '2.46'
$ws.Range('E35').Value = '  +0.25%  '
$ws.Range('E36').Value = '  -1.57%  '
Set-TextValue $ws.Range('D37') '0.852'
$ws.Range('E37').Value = '  +1.38%  '
Set-TextValue $ws.Range('D38') '0.542'
$ws.Range('E38').Value = '  -0.53%  '
$ws.Range('E40').Value = '  +1.57%  '
Set-TextValue $ws.Range('D41') '0.802'
$ws.Range('E41').Value = '  -0.75%  '
Set-TextValue $ws.Range('D42') '63.67'
$ws.Range('E42').Value = '  +1.65%  '
Set-TextValue $ws.Range('D43') '1.768.68'
$ws.Range('E43').Value = '  -1.09%  '
Set-TextValue $ws.Range('D44') '5.20'
$ws.Range('E44').Value = '  -4.05%  '
Set-TextValue $ws.Range('D45') '90.73'
$ws.Range('E45').Value = '  -1.32%  '
$ws.Range('E46').Value = '  -0.32%  '
$ws.Range('E47').Value = '  +0.44%  '
Set-TextValue $ws.Range('D48') '0.800'
$ws.Range('E48').Value = '  +19.98%  '
$ws.Range('E49').Value = '  +0.74%  '
Set-TextValue $ws.Range('D50') '7.55'
$ws.Range('E50').Value = '  -2.23%  '
Set-TextValue $ws.Range('D51') '0.0955'
$ws.Range('E51').Value = '  -2.04%  '
